$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet's tab content isn't changing, but the workbook's
# sheet name (as referenced in workbook.xml <sheet name=.../>) changes.
$ws.Name = "lake-superior-apostle-hatching"

# Update header cell A1 from "YEAR" to "year"
$ws.Range("A1").Value = "year"

# Reflect the last active selection cell as seen in the target file
$ws.Range("G31").Select()
